$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.141.40"
$ws.Range("E2").Value = "  +5.74%  "
$ws.Range("D3").Value = "1.914.53"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("D4").Formula = '="1.002"'
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Formula = '="330.64"'
$ws.Range("E5").Value = "  +5.11%  "
$ws.Range("D6").Formula = '="1.002"'
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Formula = '="0.5197"'
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("D8").Formula = '="0.4072"'
$ws.Range("E8").Value = "  +4.04%  "
$ws.Range("D9").Formula = '="0.08503"'
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").Formula = '="42.89"'
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").Formula = '="23.06"'
$ws.Range("E12").Value = "  +13.59%  "
$ws.Range("D13").Formula = '="6.449"'
$ws.Range("E13").Value = "  +4.21%  "
$ws.Range("D14").Value = "1.917.30"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Formula = '="7.378"'
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Formula = '="95.09"'
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Formula = '="0.06700"'
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Formula = '="18.41"'
$ws.Range("E20").Value = "  +4.49%  "
$ws.Range("D21").Formula = '="1.001"'
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Formula = '="6.015"'
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "30.136.54"
$ws.Range("E23").Value = "  +5.69%  "
$ws.Range("D24").Formula = '="11.34"'
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").Formula = '="2.226"'
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "2.120.61"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").Formula = '="161.67"'
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").Formula = '="21.15"'
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Formula = '="128.88"'
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("D31").Formula = '="1.091"'
$ws.Range("E31").Value = "  +5.08%  "
$ws.Range("D32").Formula = '="0.1068"'
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("D33").Formula = '="6.005"'
$ws.Range("E33").Value = "  +3.99%  "
$ws.Range("D34").Formula = '="3.616"'
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").Formula = '="0.06567"'
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Formula = '="0.2211"'
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("E38").Value = "  +3.82%  "
$ws.Range("D39").Formula = '="5.170"'
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").Formula = '="11.94"'
$ws.Range("E40").Value = "  +7.85%  "
$ws.Range("D41").Formula = '="8.806"'
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").Formula = '="0.6512"'
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("D43").Formula = '="1.240"'
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Formula = '="0.6145"'
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").Formula = '="13.27"'
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D46").Formula = '="3.749"'
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("D47").Formula = '="2.082"'
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").Formula = '="1.241"'
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Formula = '="123.89"'
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Formula = '="79.42"'
$ws.Range("E51").Value = "  +4.52%  "

# Convert the formula-forced text cells back to plain string values (strip formulas, keep text type + original style)
$ws.Range("D2:D51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
